# Natmi following Dr Hou advice
# Update ligand/receptor expression statistics in sheet1 to reflect
# the revised cell counts (1 -> 3 expressing cells) and recomputed metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.423951666666667
$ws.Range("H2").Value = 16.271855
$ws.Range("I2").Value = 0.4774188439413272
$ws.Range("J2").Value = 0.4774188439413271
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 43.69574966666666
$ws.Range("N2").Value = 131.087249
$ws.Range("O2").Value = 0.3365063034544351
$ws.Range("P2").Value = 0.3365063034544351
$ws.Range("Q2").Value = 237.0036342307661
$ws.Range("R2").Value = 2133.032708076895
$ws.Range("S2").Value = 0.1606544503741858
$ws.Range("T2").Value = 0.1606544503741858

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.423951666666667
$ws.Range("H3").Value = 16.271855
$ws.Range("I3").Value = 0.4774188439413272
$ws.Range("J3").Value = 0.4774188439413271
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("N3").Value = 140.44867
$ws.Range("O3").Value = 0.3605374521727266
$ws.Range("P3").Value = 0.3605374521727267
$ws.Range("Q3").Value = 253.9289325758722
$ws.Range("R3").Value = 2285.36039318285
$ws.Range("S3").Value = 0.1721273736138547
$ws.Range("T3").Value = 0.1721273736138547

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.423951666666667
$ws.Range("H4").Value = 16.271855
$ws.Range("I4").Value = 0.4774188439413272
$ws.Range("J4").Value = 0.4774188439413271
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.08903066666667
$ws.Range("N4").Value = 54.26709200000001
$ws.Range("O4").Value = 0.1393058338430899
$ws.Range("P4").Value = 0.1393058338430899
$ws.Range("Q4").Value = 98.11402803285112
$ws.Range("R4").Value = 883.02625229566
$ws.Range("S4").Value = 0.06650723014765059
$ws.Range("T4").Value = 0.06650723014765057

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.423951666666667
$ws.Range("H5").Value = 16.271855
$ws.Range("I5").Value = 0.4774188439413272
$ws.Range("J5").Value = 0.4774188439413271
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.25020333333334
$ws.Range("N5").Value = 63.75061
$ws.Range("O5").Value = 0.1636504105297484
$ws.Range("P5").Value = 0.1636504105297484
$ws.Range("Q5").Value = 115.2600757868389
$ws.Range("R5").Value = 1037.34068208155
$ws.Range("S5").Value = 0.07812978980563606
$ws.Range("T5").Value = 0.07812978980563606

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.583504333333333
$ws.Range("H6").Value = 4.750513
$ws.Range("I6").Value = 0.1393808158066948
$ws.Range("J6").Value = 0.1393808158066948
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 43.69574966666666
$ws.Range("N6").Value = 131.087249
$ws.Range("O6").Value = 0.3365063034544351
$ws.Range("P6").Value = 0.3365063034544351
$ws.Range("Q6").Value = 69.19240894541521
$ws.Range("R6").Value = 622.7316805087369
$ws.Range("S6").Value = 0.04690252309957436
$ws.Range("T6").Value = 0.04690252309957436

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.583504333333333
$ws.Range("H7").Value = 4.750513
$ws.Range("I7").Value = 0.1393808158066948
$ws.Range("J7").Value = 0.1393808158066948
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("N7").Value = 140.44867
$ws.Range("O7").Value = 0.3605374521727266
$ws.Range("P7").Value = 0.3605374521727267
$ws.Range("Q7").Value = 74.13369251863443
$ws.Range("R7").Value = 667.2032326677099
$ws.Range("S7").Value = 0.05025200421270185
$ws.Range("T7").Value = 0.05025200421270185

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.583504333333333
$ws.Range("H8").Value = 4.750513
$ws.Range("I8").Value = 0.1393808158066948
$ws.Range("J8").Value = 0.1393808158066948
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 18.08903066666667
$ws.Range("N8").Value = 54.26709200000001
$ws.Range("O8").Value = 0.1393058338430899
$ws.Range("P8").Value = 0.1393058338430899
$ws.Range("Q8").Value = 28.64405844646623
$ws.Range("R8").Value = 257.796526018196
$ws.Range("S8").Value = 0.01941656076768174
$ws.Range("T8").Value = 0.01941656076768174

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.583504333333333
$ws.Range("H9").Value = 4.750513
$ws.Range("I9").Value = 0.1393808158066948
$ws.Range("J9").Value = 0.1393808158066948
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 21.25020333333334
$ws.Range("N9").Value = 63.75061
$ws.Range("O9").Value = 0.1636504105297484
$ws.Range("P9").Value = 0.1636504105297484
$ws.Range("Q9").Value = 33.64978906254778
$ws.Range("R9").Value = 302.84810156293
$ws.Range("S9").Value = 0.02280972772673685
$ws.Range("T9").Value = 0.02280972772673684

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6660723333333333
$ws.Range("H10").Value = 1.998217
$ws.Range("I10").Value = 0.05862800830537802
$ws.Range("J10").Value = 0.05862800830537802
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 43.69574966666666
$ws.Range("N10").Value = 131.087249
$ws.Range("O10").Value = 0.3365063034544351
$ws.Range("P10").Value = 0.3365063034544351
$ws.Range("Q10").Value = 29.10452993722588
$ws.Range("R10").Value = 261.940769435033
$ws.Range("S10").Value = 0.01972869435373868
$ws.Range("T10").Value = 0.01972869435373868

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6660723333333333
$ws.Range("H11").Value = 1.998217
$ws.Range("I11").Value = 0.05862800830537802
$ws.Range("J11").Value = 0.05862800830537802
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.81622333333333
$ws.Range("N11").Value = 140.44867
$ws.Range("O11").Value = 0.3605374521727266
$ws.Range("P11").Value = 0.3605374521727267
$ws.Range("Q11").Value = 31.18299111348777
$ws.Range("R11").Value = 280.64692002139
$ws.Range("S11").Value = 0.02113759274038245
$ws.Range("T11").Value = 0.02113759274038245

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6660723333333333
$ws.Range("H12").Value = 1.998217
$ws.Range("I12").Value = 0.05862800830537802
$ws.Range("J12").Value = 0.05862800830537802
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.08903066666667
$ws.Range("N12").Value = 54.26709200000001
$ws.Range("O12").Value = 0.1393058338430899
$ws.Range("P12").Value = 0.1393058338430899
$ws.Range("Q12").Value = 12.04860286388489
$ws.Range("R12").Value = 108.437425774964
$ws.Range("S12").Value = 0.008167223583540283
$ws.Range("T12").Value = 0.008167223583540285

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6660723333333333
$ws.Range("H13").Value = 1.998217
$ws.Range("I13").Value = 0.05862800830537802
$ws.Range("J13").Value = 0.05862800830537802
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 21.25020333333334
$ws.Range("N13").Value = 63.75061
$ws.Range("O13").Value = 0.1636504105297484
$ws.Range("P13").Value = 0.1636504105297484
$ws.Range("Q13").Value = 14.15417251804111
$ws.Range("R13").Value = 127.38755266237
$ws.Range("S13").Value = 0.00959449762771661
$ws.Range("T13").Value = 0.009594497627716612

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.687463666666666
$ws.Range("H14").Value = 11.062391
$ws.Range("I14").Value = 0.3245723319466
$ws.Range("J14").Value = 0.3245723319466
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 43.69574966666666
$ws.Range("N14").Value = 131.087249
$ws.Range("O14").Value = 0.3365063034544351
$ws.Range("P14").Value = 0.3365063034544351
$ws.Range("Q14").Value = 161.1264892835954
$ws.Range("R14").Value = 1450.138403552359
$ws.Range("S14").Value = 0.1092206356269362
$ws.Range("T14").Value = 0.1092206356269362

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.687463666666666
$ws.Range("H15").Value = 11.062391
$ws.Range("I15").Value = 0.3245723319466
$ws.Range("J15").Value = 0.3245723319466
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 46.81622333333333
$ws.Range("N15").Value = 140.44867
$ws.Range("O15").Value = 0.3605374521727266
$ws.Range("P15").Value = 0.3605374521727267
$ws.Range("Q15").Value = 172.6331225522189
$ws.Range("R15").Value = 1553.69810296997
$ws.Range("S15").Value = 0.1170204816057876
$ws.Range("T15").Value = 0.1170204816057877

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.687463666666666
$ws.Range("H16").Value = 11.062391
$ws.Range("I16").Value = 0.3245723319466
$ws.Range("J16").Value = 0.3245723319466
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 18.08903066666667
$ws.Range("N16").Value = 54.26709200000001
$ws.Range("O16").Value = 0.1393058338430899
$ws.Range("P16").Value = 0.1393058338430899
$ws.Range("Q16").Value = 66.70264334855246
$ws.Range("R16").Value = 600.3237901369721
$ws.Range("S16").Value = 0.04521481934421726
$ws.Range("T16").Value = 0.04521481934421727

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.687463666666666
$ws.Range("H17").Value = 11.062391
$ws.Range("I17").Value = 0.3245723319466
$ws.Range("J17").Value = 0.3245723319466
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 21.25020333333334
$ws.Range("N17").Value = 63.75061
$ws.Range("O17").Value = 0.1636504105297484
$ws.Range("P17").Value = 0.1636504105297484
$ws.Range("Q17").Value = 78.35935270094556
$ws.Range("R17").Value = 705.23417430851
$ws.Range("S17").Value = 0.05311639536965884
$ws.Range("T17").Value = 0.05311639536965885
